$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet changes from the text "R40" to the text "1".
# Force the cell to keep a text (string) type -- without this, a bare numeric-
# looking string like "1" would be auto-converted to a number by Excel's
# normal cell-input parsing, which is not what the workbook change calls for
# (the stored cell keeps its shared-string / text nature).
$target = $ws.Range("B11")
$target.NumberFormat = "@"
$target.Value = "1"
